$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 31.172
$ws.Range("C2").Value = 45.144
$ws.Range("D2").Value = 76.316

$ws.Range("B3").Value = 0.981
$ws.Range("C3").Value = 1.579
$ws.Range("D3").Value = 2.56

$ws.Range("B4").Value = 6.435
$ws.Range("C4").Value = 8.468999999999999
$ws.Range("D4").Value = 14.904

$ws.Range("B5").Value = 2.56
$ws.Range("C5").Value = 3.66
$ws.Range("D5").Value = 6.220000000000001

$ws.Range("B6").Value = 41.148
$ws.Range("C6").Value = 58.852
